$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.667.20'
$ws.Range('E2').Value = '  -0.98%  '
$ws.Range('D3').Value = '3.465.40'
$ws.Range('E3').Value = '  -0.54%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '580.67'
$ws.Range('E5').Value = '  -2.13%  '
$ws.Range('D6').Value = '176.44'
$ws.Range('E6').Value = '  -1.31%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '0.598'
$ws.Range('E8').Value = '  +0.36%  '
$ws.Range('D9').Value = '3.462.65'
$ws.Range('E9').Value = '  -0.70%  '
$ws.Range('D10').Value = '0.135'
$ws.Range('E10').Value = '  -1.94%  '
$ws.Range('D11').Value = '6.88'
$ws.Range('E11').Value = '  -2.68%  '
$ws.Range('D12').Value = '0.419'
$ws.Range('E12').Value = '  -3.90%  '
$ws.Range('D13').Value = '4.071.47'
$ws.Range('E13').Value = '  -0.39%  '
$ws.Range('D14').Value = '30.35'
$ws.Range('E14').Value = '  -4.98%  '
$ws.Range('E15').Value = '  -3.00%  '
$ws.Range('D16').Value = '66.630.74'
$ws.Range('E16').Value = '  -1.11%  '
$ws.Range('D17').Value = '0.0000173'
$ws.Range('E17').Value = '  -2.42%  '
$ws.Range('D18').Value = '3.456.11'
$ws.Range('E18').Value = '  -0.68%  '
$ws.Range('D19').Value = '6.02'
$ws.Range('E19').Value = '  -3.76%  '
$ws.Range('D20').Value = '13.93'
$ws.Range('E20').Value = '  -2.94%  '
$ws.Range('D21').Value = '377.06'
$ws.Range('E21').Value = '  -3.31%  '
$ws.Range('D22').Value = '7.74'
$ws.Range('E22').Value = '  -2.95%  '
$ws.Range('E23').Value = '  +0.15%  '
$ws.Range('D24').Value = '5.73'
$ws.Range('E24').Value = '  -0.12%  '
$ws.Range('D25').Value = '71.26'
$ws.Range('E25').Value = '  -4.10%  '
$ws.Range('D26').Value = '0.528'
$ws.Range('E26').Value = '  -1.43%  '
$ws.Range('D27').Value = '0.0000118'
$ws.Range('E27').Value = '  -2.53%  '
$ws.Range('D28').Value = '9.81'
$ws.Range('E28').Value = '  -5.59%  '
$ws.Range('D29').Value = '0.172'
$ws.Range('E29').Value = '  -1.49%  '
$ws.Range('E30').Value = '  +0.09%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').Value = '24.25'
$ws.Range('E31').Value = '  +2.87%  '
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').Value = '5.88'
$ws.Range('E32').Value = '  -4.55%  '
$ws.Range('D33').Value = '1.99'
$ws.Range('E33').Value = '  -3.31%  '
$ws.Range('E34').Value = '  -5.42%  '
$ws.Range('D35').Value = '0.999'
$ws.Range('D36').Value = '7.04'
$ws.Range('E36').Value = '  -4.46%  '
$ws.Range('D37').Value = '1.52'
$ws.Range('E37').Value = '  -4.40%  '
$ws.Range('D38').Value = '158.83'
$ws.Range('E38').Value = '  -3.21%  '
$ws.Range('D39').Value = '0.877'
$ws.Range('E39').Value = '  +0.55%  '
$ws.Range('D40').Value = '27.52'
$ws.Range('E40').Value = '  +4.70%  '
$ws.Range('B41').Value = 'dogwifhat'
$ws.Range('C41').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D41').Value = '2.63'
$ws.Range('E41').Value = '  -3.67%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').Value = '1.78'
$ws.Range('E42').Value = '  -4.95%  '
$ws.Range('D43').Value = '6.54'
$ws.Range('E43').Value = '  -4.32%  '
$ws.Range('D44').Value = '4.47'
$ws.Range('E44').Value = '  -3.53%  '
$ws.Range('D45').Value = '2.697.87'
$ws.Range('E45').Value = '  -4.88%  '
$ws.Range('D46').Value = '0.0696'
$ws.Range('E46').Value = '  -3.66%  '
$ws.Range('D47').Value = '25.31'
$ws.Range('E47').Value = '  -7.16%  '
$ws.Range('D48').Value = '40.24'
$ws.Range('E48').Value = '  -3.29%  '
$ws.Range('E49').Value = '  -1.94%  '
$ws.Range('D50').Value = '323.11'
$ws.Range('E50').Value = '  -3.93%  '
$ws.Range('E51').Value = '  -3.75%  '
